# Add team record columns (Wins / Losses / Ties) to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - copy formatting from the neighboring header cell (AC1)
# then set the new header labels.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2 through 46 get the same team record values
$lastRow = 46
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 76   # AD = column 30 -> Wins
    $ws.Cells.Item($r, 31).Value = 86   # AE = column 31 -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF = column 32 -> Ties
}
